# Fix some bugs of multi-channel processing
# - change the "叠加模式" (stack/overlay mode) selector on Sheet1!B35 from "color" to "HSO"
# - add a new setting row: "通道名格式（紧邻前缀）" / "_Bin2_" (channel-name prefix format)
#   preceded by a blank separator row, matching the existing layout pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B35: overlay/stack mode color -> HSO
$ws.Range("B35").Value = "HSO"

# Row 34 (A34:B34) is an existing blank separator row with the correct label/value
# styles (A = label style, B = value style). Clone that formatting onto the two
# new rows (37 blank separator, 38 new setting) before filling in row 38's text.
$ws.Range("A34:B34").Copy()
$ws.Range("A37:B37").PasteSpecial(-4122)

$ws.Range("A34:B34").Copy()
$ws.Range("A38:B38").PasteSpecial(-4122)

$ws.Range("A38").Value = "通道名格式（紧邻前缀）"
$ws.Range("B38").Value = "_Bin2_"

# Match the final selection left behind in the saved workbook.
[void]$ws.Range("B38").Select()

Write-Host "applied multi-channel processing fixes"
